# Update cryptocurrency "Price" column (column D) values on the active sheet
# with freshly scraped figures, per the GitHub Actions symbol-list refresh.
#
# The source cells are stored as literal text (e.g. "281.80"), so each new
# value is written as Text first (NumberFormat "@") to preserve the exact
# decimal representation (trailing zeros, digit counts) instead of letting
# Excel coerce the assignment into a floating-point number. The style is
# then reset to Normal so only the cell content changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = [ordered]@{
    "D2"  = "281.74"
    "D3"  = "20.61"
    "D4"  = "6.232"
    "D5"  = "0.06145"
    "D7"  = "6.565"
    "D8"  = "1.495"
    "D9"  = "0.8188"
    "D11" = "0.1636"
    "D12" = "0.08408"
    "D13" = "0.03548"
    "D14" = "0.03190"
    "D15" = "0.09136"
    "D16" = "3.715"
    "D17" = "0.001642"
    "D18" = "0.04720"
    "D19" = "0.006525"
    "D20" = "0.006161"
    "D23" = "3.767"
    "D25" = "0.3364"
    "D40" = "0.04687"
    "D41" = "0.007190"
    "D43" = "0.1099"
    "D44" = "0.01102"
    "D45" = "0.00006551"
    "D48" = "0.002964"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}
